$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Full Timetable")

# D3: add new inline string content (was an empty inlineStr cell)
$ws.Range("D3").Value = "Cello Regulation / Maintance Workshop" + [char]10 + "(Room 201)"

# Add "Private" before "Lesson with Pitor SKWERES & pianist" and normalize double spaces
$ws.Range("C7").Value = "Venus CHAN Private Lesson with Pitor SKWERES & pianist"
$ws.Range("E7").Value = "Hannah HO Private Lesson with Pitor SKWERES & pianist"

$ws.Range("C11").Value = "Charlotte LAW Private Lesson with Pitor SKWERES & pianist"
$ws.Range("E11").Value = "Effie WONG Private Lesson with Pitor SKWERES & pianist"

$ws.Range("C19").Value = "Tucker POON Private Lesson with Pitor SKWERES & pianist"
